$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Objetivos:" (row 10) content: was wrongly holding the teacher name,
#     now holds the actual Portuguese objectives text ---
$ws.Cells.Item(10, 2).Value = "Estudo de Óptica Física."
$ws.Cells.Item(10, 3).Value = "Estudo de Óptica Física."

# --- Insert three new rows right after "Docentes responsáveis:" (row 12) to
#     hold the three professors, one per row, matching the formatting of the
#     existing "label + value" rows (columns B/C only, no A label) ---
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# Copy number/alignment formatting from the row that already has the B/C
# "value" style (row 16, the old "Programa resumido" row, now pushed down)
# into the three freshly inserted rows, then drop the stray column-A cell
# that Insert() leaves behind.
for ($r = 13; $r -le 15; $r++) {
    $ws.Cells.Item(16, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item(16, 3).Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Clear()
}

$ws.Cells.Item(13, 2).Value = "519033 - Carlos Yujiro Shigue"
$ws.Cells.Item(13, 3).Value = "519033 - Carlos Yujiro Shigue"

$ws.Cells.Item(14, 2).Value = "1341653 - Maria José Ramos Sandim"
$ws.Cells.Item(14, 3).Value = "1341653 - Maria José Ramos Sandim"

$ws.Cells.Item(15, 2).Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Cells.Item(15, 3).Value = "1643715 - Paulo Atsushi Suzuki"

# --- "Programa resumido:" (now row 16): fill in the missing Portuguese short
#     syllabus (previously incorrectly holding a professor's name) ---
$ws.Cells.Item(16, 2).Value = "Óptica de raios; Ondas eletromagnéticas: fase e polarização; Interferência; Coerência; Difração; Óptica de Fourier; Interação da luz com a matéria; Guias de ondas metálicos e dielétricos; Óptica de cristais; Óptica não linear."
$ws.Cells.Item(16, 3).Value = "Óptica de raios; Ondas eletromagnéticas: fase e polarização; Interferência; Coerência; Difração; Óptica de Fourier; Interação da luz com a matéria; Guias de ondas metálicos e dielétricos; Óptica de cristais; Óptica não linear."

# --- "Programa:" (now row 18): fill in the missing Portuguese full syllabus
#     (previously incorrectly holding a professor's name) ---
$ws.Cells.Item(18, 2).Value = "Óptica de raios. Introdução. Propagação de luz em meios homogêneos. Propagação de luz em meios não homogêneos. A lei de Snell generalizada. O princípio de Fermat. A equação dos raios. A função eikonal. Analogia ente a mecânica clássica e a óptica geométrica. O potencial óptico.Ondas eletromagnéticas. Ondas harmônicas unidimensionais. Ondas planas e esféricas. Ondas gaussianas. Propagação do feixe gaussiano. Vetor de Poynting. Intensidade.A fase da onda eletromagnética. Velocidades de fase e de grupo. Dispersão. Efeito Doppler. Aplicações astronômicas. Alargamento de linhas espectrais. Óptica relativística. Modulação eletroóptica de frequência. Automodulação de fase. Polarização das ondas eletromagnéticas. Polarização linear. Polarização elíptica. Polarização circular. Obtenção de luz linearmente polarizada. Equações de Fresnel. Polarização por reflexão total interna. Matrizes de Jones. Atividade óptica. Efeito Faraday. Isoladores ópticos. Efeito Pockels. Efeitos Kerr e Cotton-Mouton. Chaveamento eletroóptico.Interferência. Princípio da superposição. Interferência por divisão da frente de onda. Interferência por divisão de amplitudes. Interferômetro de Fabry-Perot. Analisador de espectro óptico. Teoria de películas.Coerência. Introdução. Coerência temporal. Resolução espectral de um trem de ondas finito. Coerência espacial. Medidas de diâmetros de estrelas.Difração. Princípio de Huygens. Fórmula de Fresnel-Kirchhoff. Princípio de Babinet. Difração de Fraunhofer. Difração por uma abertura circular. Rede de difração. Padrões de difração de Fresnel. Óptica de Fourier.  Microscopia por contraste de fase.  Holografia. Interação da radiação com a matéria.  Modelo do oscilador harmônico.  Dispersão cromática do índice de refração. Absorção. Espalhamento Rayleigh. Força da radiação em átomo neutro.Óptica não linear. Susceptibilidade não linear, processos paramétricos e não paramétricos. Geração de freqüências. Casamento de fase."
$ws.Cells.Item(18, 3).Value = "Óptica de raios. Introdução. Propagação de luz em meios homogêneos. Propagação de luz em meios não homogêneos. A lei de Snell generalizada. O princípio de Fermat. A equação dos raios. A função eikonal. Analogia ente a mecânica clássica e a óptica geométrica. O potencial óptico.Ondas eletromagnéticas. Ondas harmônicas unidimensionais. Ondas planas e esféricas. Ondas gaussianas. Propagação do feixe gaussiano. Vetor de Poynting. Intensidade.A fase da onda eletromagnética. Velocidades de fase e de grupo. Dispersão. Efeito Doppler. Aplicações astronômicas. Alargamento de linhas espectrais. Óptica relativística. Modulação eletroóptica de frequência. Automodulação de fase. Polarização das ondas eletromagnéticas. Polarização linear. Polarização elíptica. Polarização circular. Obtenção de luz linearmente polarizada. Equações de Fresnel. Polarização por reflexão total interna. Matrizes de Jones. Atividade óptica. Efeito Faraday. Isoladores ópticos. Efeito Pockels. Efeitos Kerr e Cotton-Mouton. Chaveamento eletroóptico.Interferência. Princípio da superposição. Interferência por divisão da frente de onda. Interferência por divisão de amplitudes. Interferômetro de Fabry-Perot. Analisador de espectro óptico. Teoria de películas.Coerência. Introdução. Coerência temporal. Resolução espectral de um trem de ondas finito. Coerência espacial. Medidas de diâmetros de estrelas.Difração. Princípio de Huygens. Fórmula de Fresnel-Kirchhoff. Princípio de Babinet. Difração de Fraunhofer. Difração por uma abertura circular. Rede de difração. Padrões de difração de Fresnel. Óptica de Fourier.  Microscopia por contraste de fase.  Holografia. Interação da radiação com a matéria.  Modelo do oscilador harmônico.  Dispersão cromática do índice de refração. Absorção. Espalhamento Rayleigh. Força da radiação em átomo neutro.Óptica não linear. Susceptibilidade não linear, processos paramétricos e não paramétricos. Geração de freqüências. Casamento de fase."

# --- "Método:" (now row 21): fix the misplaced content ---
$ws.Cells.Item(21, 2).Value = "Aulas expositivas, resolução de exercícios e seminários."
$ws.Cells.Item(21, 3).Value = "Aulas expositivas, resolução de exercícios e seminários."

# --- "Critério:" (now row 22): fix the misplaced content ---
$ws.Cells.Item(22, 2).Value = "Média ponderada de duas provas escritas: P1, P2 e TR. Conceito Final = (P1 + 2P2)/3"
$ws.Cells.Item(22, 3).Value = "Média ponderada de duas provas escritas: P1, P2 e TR. Conceito Final = (P1 + 2P2)/3"

# --- "Norma de recuperação:" (now row 23): fix the misplaced content ---
$ws.Cells.Item(23, 2).Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Cells.Item(23, 3).Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# --- "Bibliografia:" (now row 24): fill in the missing bibliography text
#     (previously incorrectly holding the recovery-exam text) ---
$ws.Cells.Item(24, 2).Value = "HECHT, E.; ZAJAC, A. Optics; Reading, Addison-Wesley, 1974.ZILLIO, S. C. Óptica Moderna - Fundamentos e Aplicações, 2005."
$ws.Cells.Item(24, 3).Value = "HECHT, E.; ZAJAC, A. Optics; Reading, Addison-Wesley, 1974.ZILLIO, S. C. Óptica Moderna - Fundamentos e Aplicações, 2005."
